$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 62926.5
$ws.Range("I53").Value = 143102.86
$ws.Range("K53").Value = 143102.86
$ws.Range("M53").Value = -142465.86

$ws.Range("H116").Value = 1773.75
$ws.Range("I116").Value = 1735
$ws.Range("J116").Value = 1786.6666
$ws.Range("K116").Value = 1735
$ws.Range("L116").Value = 1786.6666
$ws.Range("M116").Value = 1707
$ws.Range("N116").Value = -8670.6666

$ws.Range("H132").Value = 7938708
$ws.Range("I132").Value = 8930797
$ws.Range("K132").Value = 26792391
$ws.Range("M132").Value = -26789861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 43479412
$ws.Range("I2").Value = 76923896
$ws.Range("J2").Value = 1582.6
$ws.Range("K2").Value = 76923896
$ws.Range("L2").Value = 1582.6
$ws.Range("M2").Value = -76923783
$ws.Range("N2").Value = -1808.6

$ws.Range("H116").Value = 43479412
$ws.Range("I116").Value = 76923896
$ws.Range("J116").Value = 1582.6
$ws.Range("K116").Value = 76923896
$ws.Range("L116").Value = 1582.6
$ws.Range("M116").Value = -76921602
$ws.Range("N116").Value = -6170.6

$ws.Range("H132").Value = 5006.184
$ws.Range("I132").Value = 5353.069
$ws.Range("K132").Value = 16059.207
$ws.Range("M132").Value = -13529.207

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 43479412
$ws.Range("I3").Value = 76923896
$ws.Range("J3").Value = 1582.6
$ws.Range("K3").Value = 76923896
$ws.Range("L3").Value = 1582.6
$ws.Range("M3").Value = -76923782
$ws.Range("N3").Value = -1810.6

$ws.Range("H134").Value = 58120.89
$ws.Range("I134").Value = 168601.67
$ws.Range("J134").Value = 2880.5
$ws.Range("K134").Value = 505805.01
$ws.Range("L134").Value = 8641.5
$ws.Range("M134").Value = -503270.01
$ws.Range("N134").Value = -13711.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 58.166668
$ws.Range("I7").Value = 12.25
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 12.25
$ws.Range("L7").Value = 150
$ws.Range("M7").Value = 100.75
$ws.Range("N7").Value = -376

$ws.Range("H58").Value = 1282.8096
$ws.Range("I58").Value = 1291.7059
$ws.Range("J58").Value = 1245
$ws.Range("K58").Value = 1291.7059
$ws.Range("L58").Value = 1245
$ws.Range("M58").Value = -1088.7059
$ws.Range("N58").Value = -1651

$ws.Range("H107").Value = 680.7917
$ws.Range("I107").Value = 904.5
$ws.Range("J107").Value = 568.9375
$ws.Range("K107").Value = 904.5
$ws.Range("L107").Value = 568.9375
$ws.Range("M107").Value = 1015.5
$ws.Range("N107").Value = -4408.9375

$ws.Range("H132").Value = 2735.4
$ws.Range("I132").Value = 1700
$ws.Range("J132").Value = 2994.25
$ws.Range("K132").Value = 5100
$ws.Range("L132").Value = 8982.75
$ws.Range("M132").Value = -2570
$ws.Range("N132").Value = -14042.75

$ws.Range("H134").Value = 1350
$ws.Range("I134").Value = 1320
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 3960
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -1425
$ws.Range("N134").Value = -9570

$ws.Range("H136").Value = 1282.8096
$ws.Range("I136").Value = 1291.7059
$ws.Range("J136").Value = 1245
$ws.Range("K136").Value = 3875.1177
$ws.Range("L136").Value = 3735
$ws.Range("M136").Value = -1325.1177
$ws.Range("N136").Value = -8835

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 125001064
$ws.Range("I109").Value = 497.6
$ws.Range("J109").Value = 333335330
$ws.Range("K109").Value = 1492.8
$ws.Range("L109").Value = 1000005990
$ws.Range("M109").Value = -452.8000000000002
$ws.Range("N109").Value = -1000008070

$ws.Range("H131").Value = 1955.4315
$ws.Range("I131").Value = 11343.9
$ws.Range("J131").Value = 850.9059
$ws.Range("K131").Value = 34031.7
$ws.Range("L131").Value = 2552.7177
$ws.Range("M131").Value = -28991.7
$ws.Range("N131").Value = -12632.7177

$ws.Range("H132").Value = 2151
$ws.Range("I132").Value = 968.7
$ws.Range("J132").Value = 2846.4707
$ws.Range("K132").Value = 8718.300000000001
$ws.Range("L132").Value = 25618.2363
$ws.Range("M132").Value = -6188.300000000001
$ws.Range("N132").Value = -30678.2363

$ws.Range("H136").Value = 890
$ws.Range("I136").Value = 890
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2670
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 2430
$ws.Range("N136").ClearContents()

$ws.Range("H137").Value = 39751616
$ws.Range("I137").Value = 41680144
$ws.Range("J137").Value = 38564828
$ws.Range("K137").Value = 125040432
$ws.Range("L137").Value = 115694484
$ws.Range("M137").Value = -125035332
$ws.Range("N137").Value = -115704684

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 15625582
$ws.Range("I113").Value = 19231262
$ws.Range("K113").Value = 19231262
$ws.Range("M113").Value = -19229092

$ws.Range("H126").Value = 2936.6667
$ws.Range("I126").Value = 3966.6667
$ws.Range("J126").Value = 2250
$ws.Range("K126").Value = 11900.0001
$ws.Range("L126").Value = 6750
$ws.Range("M126").Value = -9430.000100000001
$ws.Range("N126").Value = -11690

$ws.Range("H132").Value = 52710.1
$ws.Range("I132").Value = 78822.30499999999
$ws.Range("K132").Value = 236466.915
$ws.Range("M132").Value = -233936.915

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1864.9166
$ws.Range("I132").Value = 1127.1765
$ws.Range("J132").Value = 3656.5715
$ws.Range("K132").Value = 3381.5295
$ws.Range("L132").Value = 10969.7145
$ws.Range("M132").Value = -851.5295000000001
$ws.Range("N132").Value = -16029.7145

$ws.Range("H133").Value = 24756.5
$ws.Range("J133").Value = 24756.5
$ws.Range("L133").Value = 24756.5
$ws.Range("N133").Value = -29816.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3742.0977
$ws.Range("I136").Value = 4329.273
$ws.Range("J136").Value = 1320
$ws.Range("K136").Value = 12987.819
$ws.Range("L136").Value = 3960
$ws.Range("M136").Value = -10437.819
$ws.Range("N136").Value = -9060

$ws.Range("H140").Value = 52865.6
$ws.Range("J140").Value = 52865.6
$ws.Range("L140").Value = 52865.6
$ws.Range("N140").Value = -63225.6
